$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 328.42856
$ws.Range("I32").Value = 310
$ws.Range("J32").Value = 374.5
$ws.Range("K32").Value = 310
$ws.Range("L32").Value = 374.5
$ws.Range("M32").Value = -1026.5
$ws.Range("H88").Value = 6024.28
$ws.Range("I88").Value = 5369.364
$ws.Range("J88").Value = 6538.857
$ws.Range("K88").Value = 5369.364
$ws.Range("L88").Value = 6538.857
$ws.Range("M88").Value = -4963.364
$ws.Range("N88").Value = -7350.857
$ws.Range("H91").Value = 6024.28
$ws.Range("I91").Value = 5369.364
$ws.Range("J91").Value = 6538.857
$ws.Range("K91").Value = 5369.364
$ws.Range("L91").Value = 6538.857
$ws.Range("M91").Value = -3965.364
$ws.Range("N91").Value = -9346.857
$ws.Range("H96").Value = 802.2941
$ws.Range("I96").Value = 758.1
$ws.Range("J96").Value = 865.4286
$ws.Range("K96").Value = 2274.3
$ws.Range("L96").Value = 2596.2858
$ws.Range("M96").Value = -901.3000000000002
$ws.Range("N96").Value = -5342.2858
$ws.Range("H100").Value = 3150.8823
$ws.Range("I100").Value = 2400.3572
$ws.Range("J100").Value = 6653.3335
$ws.Range("K100").Value = 2400.3572
$ws.Range("L100").Value = 6653.3335
$ws.Range("M100").Value = -1859.3572
$ws.Range("N100").Value = -7735.3335
$ws.Range("H132").Value = 3280531.5
$ws.Range("I132").Value = 3449995.2
$ws.Range("K132").Value = 10349985.6
$ws.Range("M132").Value = -10347455.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5263.6123
$ws.Range("I32").Value = 3863.965
$ws.Range("J32").Value = 15294.417
$ws.Range("K32").Value = 3863.965
$ws.Range("L32").Value = 15294.417
$ws.Range("M32").Value = -3576.965
$ws.Range("N32").Value = -15868.417
$ws.Range("H45").Value = 1301.3846
$ws.Range("I45").Value = 1015.80646
$ws.Range("J45").Value = 2408
$ws.Range("K45").Value = 1015.80646
$ws.Range("L45").Value = 2408
$ws.Range("M45").Value = -638.80646
$ws.Range("N45").Value = -3162
$ws.Range("H61").Value = 4847.0586
$ws.Range("I61").Value = 2200
$ws.Range("J61").Value = 5200
$ws.Range("K61").Value = 2200
$ws.Range("L61").Value = 5200
$ws.Range("M61").Value = -1988
$ws.Range("N61").Value = -5624
$ws.Range("H74").Value = 848
$ws.Range("I74").Value = 775.5333000000001
$ws.Range("J74").Value = 1029.1666
$ws.Range("K74").Value = 775.5333000000001
$ws.Range("L74").Value = 1029.1666
$ws.Range("M74").Value = 98.46669999999995
$ws.Range("N74").Value = -2777.1666
$ws.Range("H77").Value = 848
$ws.Range("I77").Value = 775.5333000000001
$ws.Range("J77").Value = 1029.1666
$ws.Range("K77").Value = 3877.6665
$ws.Range("L77").Value = 5145.833000000001
$ws.Range("M77").Value = 490.3334999999997
$ws.Range("N77").Value = -13881.833
$ws.Range("H97").Value = 570.2593000000001
$ws.Range("I97").Value = 564.85
$ws.Range("J97").Value = 585.7143
$ws.Range("K97").Value = 564.85
$ws.Range("L97").Value = 585.7143
$ws.Range("M97").Value = -68.85000000000002
$ws.Range("N97").Value = -1577.7143
$ws.Range("H122").Value = 3362.7144
$ws.Range("I122").Value = 2175.3333
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 6525.999899999999
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -4075.999899999999
$ws.Range("N122").Value = -21400
$ws.Range("H136").Value = 4847.0586
$ws.Range("I136").Value = 2200
$ws.Range("J136").Value = 5200
$ws.Range("K136").Value = 6600
$ws.Range("L136").Value = 15600
$ws.Range("M136").Value = -4050
$ws.Range("N136").Value = -20700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 196.875
$ws.Range("I22").Value = 95.833336
$ws.Range("K22").Value = 95.833336
$ws.Range("M22").Value = 77.166664
$ws.Range("H45").Value = 30065
$ws.Range("J45").Value = 30065
$ws.Range("L45").Value = 30065
$ws.Range("N45").Value = -31681
$ws.Range("H86").Value = 15374.595
$ws.Range("I86").Value = 1232.4348
$ws.Range("J86").Value = 38608.145
$ws.Range("K86").Value = 1232.4348
$ws.Range("L86").Value = 38608.145
$ws.Range("M86").Value = -109.4348
$ws.Range("N86").Value = -40854.145
$ws.Range("H89").Value = 15374.595
$ws.Range("I89").Value = 1232.4348
$ws.Range("J89").Value = 38608.145
$ws.Range("K89").Value = 6162.174
$ws.Range("L89").Value = 193040.725
$ws.Range("M89").Value = -546.174
$ws.Range("N89").Value = -204272.725
$ws.Range("H94").Value = 660.9643
$ws.Range("I94").Value = 578.6111
$ws.Range("J94").Value = 809.2
$ws.Range("K94").Value = 578.6111
$ws.Range("L94").Value = 809.2
$ws.Range("M94").Value = -127.6111
$ws.Range("N94").Value = -1711.2
$ws.Range("H99").Value = 2962.5
$ws.Range("I99").Value = 2283.3333
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 2283.3333
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -785.3332999999998
$ws.Range("N99").Value = -7996
$ws.Range("H134").Value = 3083.5264
$ws.Range("I134").Value = 2113.4285
$ws.Range("J134").Value = 5799.8
$ws.Range("K134").Value = 6340.2855
$ws.Range("L134").Value = 17399.4
$ws.Range("M134").Value = -3805.2855
$ws.Range("N134").Value = -22469.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2285.077
$ws.Range("I75").Value = 1012.75
$ws.Range("J75").Value = 2850.5557
$ws.Range("K75").Value = 3038.25
$ws.Range("L75").Value = 8551.667099999999
$ws.Range("M75").Value = -2040.25
$ws.Range("N75").Value = -10547.6671
$ws.Range("H78").Value = 2285.077
$ws.Range("I78").Value = 1012.75
$ws.Range("J78").Value = 2850.5557
$ws.Range("K78").Value = 9114.75
$ws.Range("L78").Value = 25655.0013
$ws.Range("M78").Value = -4122.75
$ws.Range("N78").Value = -35639.0013
$ws.Range("H120").Value = 19505
$ws.Range("I120").Value = 19010
$ws.Range("K120").Value = 57030
$ws.Range("M120").Value = -52192

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4260
$ws.Range("I80").Value = 4016
$ws.Range("J80").Value = 4666.6665
$ws.Range("K80").Value = 4016
$ws.Range("L80").Value = 4666.6665
$ws.Range("M80").Value = -3018
$ws.Range("N80").Value = -6662.6665
$ws.Range("H83").Value = 4260
$ws.Range("I83").Value = 4016
$ws.Range("J83").Value = 4666.6665
$ws.Range("K83").Value = 20080
$ws.Range("L83").Value = 23333.3325
$ws.Range("M83").Value = -15088
$ws.Range("N83").Value = -33317.3325
$ws.Range("H122").Value = 7172.222
$ws.Range("I122").Value = 4275
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 12825
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -10375
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 3623.4856
$ws.Range("I132").Value = 3286.8572
$ws.Range("J132").Value = 4128.4287
$ws.Range("K132").Value = 9860.571599999999
$ws.Range("L132").Value = 12385.2861
$ws.Range("M132").Value = -7330.571599999999
$ws.Range("N132").Value = -17445.2861
$ws.Range("H137").Value = 29642.857
$ws.Range("J137").Value = 29642.857
$ws.Range("L137").Value = 29642.857
$ws.Range("N137").Value = -39842.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3100
$ws.Range("I40").Value = 1400
$ws.Range("J40").Value = 3440
$ws.Range("K40").Value = 1400
$ws.Range("L40").Value = 3440
$ws.Range("M40").Value = -1264
$ws.Range("N40").Value = -3712
$ws.Range("H61").Value = 58826876
$ws.Range("I61").Value = 111112790
$ws.Range("K61").Value = 111112790
$ws.Range("M61").Value = -111112588
$ws.Range("H113").Value = 58826876
$ws.Range("I113").Value = 111112790
$ws.Range("K113").Value = 111112790
$ws.Range("M113").Value = -111110620
$ws.Range("H122").Value = 3091.7778
$ws.Range("I122").Value = 2592.4614
$ws.Range("J122").Value = 4390
$ws.Range("K122").Value = 7777.3842
$ws.Range("L122").Value = 13170
$ws.Range("M122").Value = -5327.3842
$ws.Range("N122").Value = -18070

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 70005
$ws.Range("J14").Value = 70005
$ws.Range("L14").Value = 70005
$ws.Range("N14").Value = -70341
$ws.Range("H132").Value = 15968.025
$ws.Range("I132").Value = 3352.5186
$ws.Range("J132").Value = 42169.46
$ws.Range("K132").Value = 10057.5558
$ws.Range("L132").Value = 126508.38
$ws.Range("M132").Value = -7527.5558
$ws.Range("N132").Value = -131568.38
$ws.Range("H133").Value = 39263.332
$ws.Range("J133").Value = 39263.332
$ws.Range("L133").Value = 39263.332
$ws.Range("N133").Value = -49383.332
$ws.Range("H135").Value = 76925.71000000001
$ws.Range("J135").Value = 76925.71000000001
$ws.Range("L135").Value = 76925.71000000001
$ws.Range("N135").Value = -87065.71000000001
